$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "data_source": insert a new row for "Wind Speed at 100m"
# (Global Wind Atlas) above the existing "Water Area" row, and add
# an "OSM" / geofabrik link source to the five rows that follow.
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data_source")

# Push rows 5-9 down to 6-10 to make room for the new layer entry.
$wsData.Rows("5:5").Insert()

# New row 5: Wind Speed at 100m / Global Wind Atlas / link
# (Order of assignment matters so shared-string table ends up in
# the same order as produced by a human editing the sheet: the
# label, then the link, then the source name.)
$wsData.Range("A5").Value = "Wind Speed at 100m"
$wsData.Range("C5").Value = "https://globalwindatlas.info/en"
$wsData.Range("B5").Value = "Global Wind Atlas"

# Give the new link cell the same "hyperlink look" style already
# used by the other Link cells (C2:C4) by copying their format.
$wsData.Range("C2").Copy()
$wsData.Range("C5").PasteSpecial(-4122)

# Rows 6-10 (previously 5-9) now also reference OSM / geofabrik as
# their data source.
$wsData.Range("B6").Value = "OSM"
$wsData.Range("C6").Value = "https://download.geofabrik.de/europe/germany/hamburg.html"
$wsData.Range("B7").Value = "OSM"
$wsData.Range("C7").Value = "https://download.geofabrik.de/europe/germany/hamburg.html"
$wsData.Range("B8").Value = "OSM"
$wsData.Range("C8").Value = "https://download.geofabrik.de/europe/germany/hamburg.html"
$wsData.Range("B9").Value = "OSM"
$wsData.Range("C9").Value = "https://download.geofabrik.de/europe/germany/hamburg.html"
$wsData.Range("B10").Value = "OSM"
$wsData.Range("C10").Value = "https://download.geofabrik.de/europe/germany/hamburg.html"

$wsData.Range("C2").Copy()
$wsData.Range("C6:C10").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet "MCA_criteria": add the matching MCA layer row for the new
# Wind Speed at 100m data source.
# ---------------------------------------------------------------
$wsMca = $wb.Worksheets.Item("MCA_criteria")
$wsMca.Range("A9").Value = "Wind Speed at 100m"
$wsMca.Range("C9").Value = "Raster"

# Column A needs to widen to fit the new, longer label.
$wsMca.Columns("A:A").ColumnWidth = 18.79

# ---------------------------------------------------------------
# Restore the last-used selection on each sheet and make sure the
# data_source tab is the active one, as in the source workbook.
# ---------------------------------------------------------------
$null = $wsMca.Range("B4").Select()
$null = $wsData.Activate()
$null = $wsData.Range("B13").Select()
